$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '305.66'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.86%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '38.83'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '6.70%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.113'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.02%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08066'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.16%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.933'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-2.27%'
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.027'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2.41%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9268'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.19%'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1444'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-2.56%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1914'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-1.12%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09039'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.48%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03506'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.69%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09774'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.90%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001401'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-1.02%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005897'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-3.03%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.758'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-2.08%'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.208'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.06%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.378'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.66%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3462'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.29%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1327'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.03%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.695'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.46%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04377'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.60%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001207'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-2.21%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004267'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '2.45%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001302'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.00%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02033'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-1.60%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05046'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-1.25%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007517'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '0.37%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009777'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-3.36%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1340'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-1.83%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002094'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-1.41%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.009832'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '1.63%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006211'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000751'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.01%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002873'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '12.52%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002103'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.01%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002003'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.01%'
